# Apply data corrections to Tab17 as captured by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minor floating point correction on row 67 (I67)
$ws.Range("I67").Value = 75378.546306000004

# Row 97 ("Afrique, États fragiles") - corrected aggregate values
$ws.Range("C97").Value = 144193.31095099999
$ws.Range("D97").Value = 78305.983385
$ws.Range("E97").Value = 22463.926414000001
$ws.Range("F97").Value = 244963.22075000001
$ws.Range("G97").Value = 19017.378799999999
$ws.Range("H97").Value = 159093.021412
$ws.Range("I97").Value = 94732.637432999996
$ws.Range("J97").Value = 272843.03764499997

# Row 98 ("RDM, États fragiles") - corrected aggregate values
$ws.Range("C98").Value = 40674.361312000001
$ws.Range("D98").Value = 59146.964883000001
$ws.Range("E98").Value = 112031.16651700001
$ws.Range("F98").Value = 211852.49271200001
$ws.Range("G98").Value = 39684.129846999997
$ws.Range("H98").Value = 203628.562393
$ws.Range("I98").Value = 84104.150620999993
$ws.Range("J98").Value = 327416.84286099998
